$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(18, 9).Value = "aa"
$ws.Cells.Item(18, 10).Value = "Agree/Accept"
$ws.Cells.Item(20, 9).Value = "aa"
$ws.Cells.Item(20, 10).Value = "Agree/Accept"
$ws.Cells.Item(21, 9).Value = "aa"
$ws.Cells.Item(21, 10).Value = "Agree/Accept"
$ws.Cells.Item(22, 9).Value = "sd"
$ws.Cells.Item(22, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(34, 9).Value = "sv"
$ws.Cells.Item(34, 10).Value = "Statement-opinion"
$ws.Cells.Item(47, 9).Value = "sd"
$ws.Cells.Item(47, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(50, 9).Value = "sd"
$ws.Cells.Item(50, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(52, 9).Value = "aa"
$ws.Cells.Item(52, 10).Value = "Agree/Accept"
$ws.Cells.Item(69, 9).Value = "%"
$ws.Cells.Item(69, 10).Value = "Uninterpretable"
$ws.Cells.Item(77, 9).Value = "sd"
$ws.Cells.Item(77, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(112, 9).Value = "sv"
$ws.Cells.Item(112, 10).Value = "Statement-opinion"
$ws.Cells.Item(132, 9).Value = "sd"
$ws.Cells.Item(132, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(186, 9).Value = "ba"
$ws.Cells.Item(186, 10).Value = "Appreciation"
$ws.Cells.Item(187, 9).Value = "sd"
$ws.Cells.Item(187, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(194, 9).Value = "ba"
$ws.Cells.Item(194, 10).Value = "Appreciation"
$ws.Cells.Item(219, 9).Value = "sd"
$ws.Cells.Item(219, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(224, 9).Value = "sd"
$ws.Cells.Item(224, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(231, 9).Value = "sd"
$ws.Cells.Item(231, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(236, 9).Value = "aa"
$ws.Cells.Item(236, 10).Value = "Agree/Accept"
$ws.Cells.Item(237, 9).Value = "aa"
$ws.Cells.Item(237, 10).Value = "Agree/Accept"
$ws.Cells.Item(241, 9).Value = "sd"
$ws.Cells.Item(241, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(265, 9).Value = "sv"
$ws.Cells.Item(265, 10).Value = "Statement-opinion"
$ws.Cells.Item(287, 9).Value = "aa"
$ws.Cells.Item(287, 10).Value = "Agree/Accept"
$ws.Cells.Item(298, 9).Value = "sv"
$ws.Cells.Item(298, 10).Value = "Statement-opinion"
$ws.Cells.Item(300, 9).Value = "sd"
$ws.Cells.Item(300, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(302, 9).Value = "ba"
$ws.Cells.Item(302, 10).Value = "Appreciation"
$ws.Cells.Item(303, 9).Value = "sd"
$ws.Cells.Item(303, 10).Value = "Statement-non-opinion"